$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph "Входные данные – статья"  ->  "Входные данные – статья, автор."
# Paragraph "Выходные данные – статья, номер." -> split into its own two
# paragraphs:
#   "Входные данные – статья, автор."
#   "Выходные данные – кол-во статей, сумма выплаченная автору."
# ---------------------------------------------------------------------------

# Locate the two paragraphs by their current text so the script is resilient
# to exact indices.
$targetIn = $null
$targetOut = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Входные данные")) { $targetIn = $p }
    if ($t.StartsWith("Выходные данные")) { $targetOut = $p }
}

$pOut = $targetOut
$pOutRange = $pOut.Range
$pOutStart = $pOutRange.Start

# --- 1) Split "статья" out of "Выходные данные – статья" into its own run,
#        replacing it with "кол-во статей". Bracketing the edit with a
#        temporary bookmark stops the engine from re-coalescing the new run
#        back into the preceding "Выходные данные – " run (they'd otherwise
#        merge because both end up with identical/default formatting).
$searchRange = $d.Range($pOutStart, $pOutRange.End)
$found = $searchRange.Find.Execute("статья", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$wordStart = $searchRange.Start
$wordEnd = $searchRange.End
$target = $d.Range($wordStart, $wordEnd)
$target.Text = "кол-во статей"

$splitPoint = $d.Range($wordStart, $wordStart)
$d.Bookmarks.Add("zzz_tmp_split", $splitPoint)
$d.Bookmarks("zzz_tmp_split").Delete()

# --- 2) " номер." -> " сумма выплаченная автору" (drop the trailing period;
#        it gets re-added after the _GoBack bookmark below).
$pOutRange = $pOut.Range
$searchRange2 = $d.Range($pOutRange.Start, $pOutRange.End)
$searchRange2.Find.Execute(" номер.", $true, $false, $false, $false, $false, $true, 1, $false, " сумма выплаченная автору", 2)

# --- 3) Re-add the trailing "." after the (now relocated) _GoBack bookmark,
#        at the very end of the paragraph.
$pOutRange = $pOut.Range
$insPos = $pOutRange.End - 1
$periodTarget = $d.Range($insPos, $insPos)
$periodTarget.InsertAfter(".")

# --- 4) Append ", автор." to the "Входные данные – статья" paragraph, using
#        the comma's FormattedText from the (already en-US tagged) comma run
#        in the "Выходные" paragraph so the copy keeps that exact formatting.
$pOutRange = $pOut.Range
$commaSearch = $d.Range($pOutRange.Start, $pOutRange.End)
$commaSearch.Find.Execute(",", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$commaFormatted = $commaSearch.FormattedText

$pIn = $targetIn
$pInRange = $pIn.Range
$insPos2 = $pInRange.End - 1
$commaTarget = $d.Range($insPos2, $insPos2)
$commaTarget.FormattedText = $commaFormatted

$pInRange = $pIn.Range
$insPos3 = $pInRange.End - 1
$textTarget = $d.Range($insPos3, $insPos3)
$textTarget.InsertAfter(" автор.")

Write-Host "edit complete"
